$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 239, shifting existing rows 239-333 down to 240-334
$ws.Rows("239:239").Insert()

# Populate the newly inserted row 239 with the new data record
$ws.Range("A239").Value = 4
$ws.Range("B239").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C239").Value = "Los Lagos"
$ws.Range("D239").Value = 45119
$ws.Range("E239").Value = 10
$ws.Range("F239").Value = 100112009
$ws.Range("G239").Value = "Acelga"
$ws.Range("H239").Value = "Sin especificar"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 15
$ws.Range("K239").Value = 10000
$ws.Range("L239").Value = 10000
$ws.Range("M239").Value = 10000
$ws.Range("N239").Value = "$/docena de atados (12 kilos)"
$ws.Range("O239").Value = "Región de La Araucanía"
$ws.Range("P239").Value = 833
$ws.Range("Q239").Value = 12
$ws.Range("R239").Value = "Hortaliza"

# Apply the same date style used in column D for the other date cells (style index 2)
$ws.Range("D239").NumberFormat = $ws.Range("D240").NumberFormat
